$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: new shared-strings are appended to the sharedStrings table in the
# order their cells are first written, so the write order below is chosen
# to reproduce the same table ordering as the target workbook (trustbot1.0,
# mybot(2) v3.4.1, 125games, d7, d8, New D5 battles:, 3.4.1, 5s+40ms,
# Sample size, 5+0.05s, tied).

# --- Row 45 first (introduces "trustbot1.0") ---
$ws.Range("A45").Value = "trustbot1.0"
$ws.Range("B45").Formula = "=14-95"

# --- Row 44 (introduces "mybot(2) v3.4.1" then "125games") ---
$ws.Range("B44").Value = "mybot(2) v3.4.1"
$ws.Range("A44").Value = "125games"

# --- Row 46 ---
$ws.Range("A46").Value = 1.1
$ws.Range("B46").Formula = "=25-74"

# --- Row 47 ---
$ws.Range("A47").Value = 1.2
$ws.Range("B47").Formula = "=32-77"

# --- Row 48 ---
$ws.Range("A48").Value = 1.3
$ws.Range("B48").Formula = "=33-76"

# --- Row 49 ---
$ws.Range("A49").Value = 1.4
$ws.Range("B49").Formula = "=(56+16)/182"

# --- Row 50 ---
$ws.Range("A50").Value = 1.5
$ws.Range("B50").Formula = "=(73+12)/216"

# --- Row 51 ---
$ws.Range("A51").Value = 1.6
$ws.Range("B51").Formula = "=(85+14)/201"

# --- Row 52 ---
$ws.Range("A52").Value = 1.7
$ws.Range("B52").Formula = "=(87+19)/206"

# --- Row 53 ---
$ws.Range("A53").Value = 1.8
$ws.Range("B53").Formula = "=(222+95/2)/560"

# --- Row 54 ---
$ws.Range("A54").Value = 1.9
$ws.Range("B54").Formula = "=(397+93)/1000"

# --- Row 55 ---
$ws.Range("A55").Value = 2
$ws.Range("B55").Value = 0.5

# --- Row 57 ---
$ws.Range("B57").Value = 2

# --- Row 58 ---
$ws.Range("A58").Value = 2.1

# --- Row 60 (introduces "d7") ---
$ws.Range("A60").Value = "d7"
$ws.Range("B60").Value = 2

# --- Row 61 ---
$ws.Range("A61").Value = 2.1
$ws.Range("B61").Value = 0.49

# --- Row 62 ---
$ws.Range("A62").Value = 2.2
$ws.Range("B62").Value = 0.49

# --- Row 64 (introduces "d8") ---
$ws.Range("A64").Value = "d8"

# --- Row 65 ---
$ws.Range("A65").Value = 1.8
$ws.Range("B65").Formula = "=(315+157/2)/755"

# --- Row 67: header row for "New D5 battles:" table (introduces the
#     remaining header strings, in left-to-right order) ---
$ws.Range("A67").Value = "New D5 battles:"
$ws.Range("B67").Value = "3.4.1"
$ws.Range("C67").HorizontalAlignment = -4152
$ws.Range("C67").Value = "5s+40ms"
$ws.Range("D67").Value = "Sample size"

# --- Row 68 ---
$ws.Range("A68").Value = 2.1
$ws.Range("B68").Formula = "=(112+22)/252"
$ws.Range("C68").HorizontalAlignment = -4152
$ws.Range("C68").Formula = "=(95+25)/312"
$ws.Range("D68").Value = 312

# --- Row 69 ---
$ws.Range("A69").Value = 0.1
$ws.Range("C69").Formula = "=(117+32)/329"
$ws.Range("D69").Value = 329

# --- Row 70 ---
$ws.Range("A70").Value = 0.2
$ws.Range("C70").Formula = "=(101+24)/250"
$ws.Range("D70").Value = 250

# --- Row 71 ---
$ws.Range("A71").Value = 0.3
$ws.Range("C71").Formula = "=(159+30)/369"

# --- Row 73 (introduces "5+0.05s") ---
$ws.Range("A73").HorizontalAlignment = -4152
$ws.Range("A73").Value = "5+0.05s"
$ws.Range("B73").Value = 0.3

# --- Row 74 (introduces "tied") ---
$ws.Range("A74").Value = 0.4
$ws.Range("B74").Value = "tied"
